$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" data column (Q) to the right of the existing "2019" column (P),
# matching formatting of the corresponding P-column cell in each row.
$ws.Range("P4:P8").Copy()
$ws.Range("Q4:Q8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 53.2
$ws.Range("Q6").Value = 23.2
$ws.Range("Q7").Value = 10
$ws.Range("Q8").Value = 20

# Rows 7 and 8 display their figures with one decimal place (matching the "0.0"
# number format already used by the rest of the data rows in this table).
$ws.Range("Q7").NumberFormat = "0.0"
$ws.Range("Q8").NumberFormat = "0.0"

# Match the author's final selection.
$ws.Range("P9").Select() | Out-Null
